{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Changes implemented (per the diff):\n//  1. Remove the \"_GoBack\" bookmark from the very first paragraph (the\n//     document title \"\u5317\u4eac\u5de5\u4e1a\u5927\u5b66\").\n//  2. Delete the empty paragraph (style \"List Paragraph\") that sits\n//     between the \"...\u68c0\u6d4b\u65b9\u6cd5\u8bbe\u8ba1\" paragraph and the \"\u901a\u8fc7\u5bf9ROP\u653b\u51fb...\"\n//     paragraph.\n//  3. Re-insert the \"_GoBack\" bookmark at the start of the paragraph that\n//     begins with \"\u5b9e\u73b0\u4e00\u4e2a\u68c0\u6d4bROP\u653b\u51fb\u3001...\" (the last edited spot).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Step 1: delete the stray \"_GoBack\" bookmark wherever it currently lives.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Step 2: find & delete the empty \"List Paragraph\" that follows the\n// paragraph ending in \"\u68c0\u6d4b\u65b9\u6cd5\u8bbe\u8ba1\" and precedes \"\u901a\u8fc7\u5bf9ROP\u653b\u51fb...\".\nlet emptyParaIndex = -1;\nfor (let i = 0; i < items.length - 1; i++) {\n  if (\n    items[i].text.indexOf(\"\u68c0\u6d4b\u65b9\u6cd5\u8bbe\u8ba1\") !== -1 &&\n    items[i + 1].text === \"\"\n  ) {\n    emptyParaIndex = i + 1;\n    break;\n  }\n}\nif (emptyParaIndex !== -1) {\n  items[emptyParaIndex].delete();\n  await context.sync();\n}\n\n// Step 3: add the \"_GoBack\" bookmark back at the start of the paragraph\n// that begins with \"\u5b9e\u73b0\u4e00\u4e2a\u68c0\u6d4bROP\u653b\u51fb\u3001\".\nlet targetParagraph = null;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"\u5b9e\u73b0\u4e00\u4e2a\u68c0\u6d4bROP\u653b\u51fb\") === 0) {\n    targetParagraph = items[i];\n    break;\n  }\n}\nif (targetParagraph) {\n  const startRange = targetParagraph.getRange(\"Start\");\n  startRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n#\n# Changes implemented (per the diff):\n#  1. Remove the \"_GoBack\" bookmark from the very first paragraph (the\n#     document title \"\u5317\u4eac\u5de5\u4e1a\u5927\u5b66\").\n#  2. Delete the empty paragraph (style \"List Paragraph\") that sits\n#     between the \"...\u68c0\u6d4b\u65b9\u6cd5\u8bbe\u8ba1\" paragraph and the \"\u901a\u8fc7\u5bf9ROP\u653b\u51fb...\"\n#     paragraph.\n#  3. Re-insert the \"_GoBack\" bookmark at the start of the paragraph that\n#     begins with \"\u5b9e\u73b0\u4e00\u4e2a\u68c0\u6d4bROP\u653b\u51fb\u3001...\" (the last edited spot).\n\n$d = $word.ActiveDocument\n\n# Step 1: delete the stray \"_GoBack\" bookmark wherever it currently lives.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Step 2: find & delete the empty paragraph that follows the paragraph\n# ending in \"\u68c0\u6d4b\u65b9\u6cd5\u8bbe\u8ba1\" and precedes \"\u901a\u8fc7\u5bf9ROP\u653b\u51fb...\".\n$paras = $d.Paragraphs\n$prevHadMarker = $false\n$emptyPara = $null\nforeach ($p in $paras) {\n    if ($prevHadMarker -and $p.Range.Text.Length -le 1) {\n        $emptyPara = $p\n        break\n    }\n    $prevHadMarker = $p.Range.Text.Contains(\"\u68c0\u6d4b\u65b9\u6cd5\u8bbe\u8ba1\")\n}\nif ($emptyPara -ne $null) {\n    $emptyPara.Range.Delete()\n}\n\n# Step 3: add the \"_GoBack\" bookmark back at the start of the paragraph\n# that begins with \"\u5b9e\u73b0\u4e00\u4e2a\u68c0\u6d4bROP\u653b\u51fb\u3001\".\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"\u5b9e\u73b0\u4e00\u4e2a\u68c0\u6d4bROP\u653b\u51fb\")) {\n        $target = $p\n        break\n    }\n}\nif ($target -ne $null) {\n    $startRange = $d.Range($target.Range.Start, $target.Range.Start)\n    $d.Bookmarks.Add(\"_GoBack\", $startRange)\n}\n"}
